# Apply updated crypto price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.162.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "'2.449.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'583.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "'142.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'2.442.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "'0.343"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "'0.0000176"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'2.878.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'62.047.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "'2.435.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  -4.08%  "
$ws.Range("D20").Value = "'7.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "'326.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -6.48%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "'9.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").Value = "'599.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.88%  "
$ws.Range("D28").Value = "'0.0₃0968"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").Value = "'7.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").Value = "'1.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").Value = "'4.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("D39").Value = "'152.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.54%  "
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").Value = "'5.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").Value = "'43.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").Value = "'0.0₆0277"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +24.02%  "
$ws.Range("D47").Value = "'141.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").Value = "'0.0519"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").Value = "'19.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.77%  "
